# Update crypto price/volume data (and two row re-orderings for
# EnergySwap/Decentraland and EOS/NEARProtocol) per the Apr 6 2023
# GitHub Actions data refresh.
#
# Numeric-looking Price values are written with a leading apostrophe
# (forcing text entry, like a user typing it in Excel) so they stay
# text instead of being auto-coerced to numbers; the style is then
# reset to Normal so no stray quote-prefix formatting is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.946.89'
$ws.Range("E2").Value = '  -1.25%  '
$ws.Range("D3").Value = '1.870.09'
$ws.Range("E3").Value = '  -1.90%  '
$ws.Range("E4").Value = '  -0.24%  '
$ws.Range("D5").Value = '''312.61'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.74%  '
$ws.Range("E6").Value = '  -0.24%  '
$ws.Range("D7").Value = '''0.5034'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.97%  '
$ws.Range("E8").Value = '  -2.81%  '
$ws.Range("D9").Value = '''0.08925'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -7.77%  '
$ws.Range("D10").Value = '''1.117'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.05%  '
$ws.Range("D11").Value = '''41.55'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.31%  '
$ws.Range("D12").Value = '''6.385'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.78%  '
$ws.Range("D13").Value = '''20.67'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.32%  '
$ws.Range("D14").Value = '1.862.73'
$ws.Range("E14").Value = '  -3.43%  '
$ws.Range("D15").Value = '''7.240'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.42%  '
$ws.Range("D16").Value = '''1.000'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.24%  '
$ws.Range("D17").Value = '''0.00001099'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.23%  '
$ws.Range("D18").Value = '''91.15'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.94%  '
$ws.Range("D19").Value = '''0.06659'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.34%  '
$ws.Range("D20").Value = '''18.13'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.35%  '
$ws.Range("D21").Value = '''1.000'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.28%  '
$ws.Range("D22").Value = '''6.116'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.81%  '
$ws.Range("D23").Value = '27.963.02'
$ws.Range("E23").Value = '  -1.40%  '
$ws.Range("D24").Value = '''11.48'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.27%  '
$ws.Range("D25").Value = '''2.279'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.56%  '
$ws.Range("D26").Value = '2.075.81'
$ws.Range("E26").Value = '  -2.80%  '
$ws.Range("D27").Value = '''2.495'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -6.59%  '
$ws.Range("D28").Value = '''158.38'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.25%  '
$ws.Range("D29").Value = '''20.67'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.63%  '
$ws.Range("D30").Value = '''126.23'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.02%  '
$ws.Range("D31").Value = '''0.1061'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.40%  '
$ws.Range("D32").Value = '''1.056'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.85%  '
$ws.Range("D33").Value = '''5.604'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.92%  '
$ws.Range("D34").Value = '''3.608'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.73%  '
$ws.Range("D35").Value = '''9.514'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.01%  '
$ws.Range("D36").Value = '''0.06556'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.90%  '
$ws.Range("D37").Value = '''0.02397'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.53%  '
$ws.Range("D38").Value = '''0.2182'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.23%  '
$ws.Range("D39").Value = '''1.282'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.23%  '
$ws.Range("E40").Value = '  -3.45%  '
$ws.Range("D41").Value = '''0.6372'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.48%  '
$ws.Range("D42").Value = '''11.50'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.30%  '
$ws.Range("D43").Value = '''4.902'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.53%  '
$ws.Range("D44").Value = '''0.9999'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.36%  '
$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").Value = '''13.17'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.60%  '
$ws.Range("B46").Value = 'Decentraland'
$ws.Range("C46").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D46").Value = '''0.6004'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.41%  '
$ws.Range("D47").Value = '''1.279'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.16%  '
$ws.Range("D48").Value = '''3.666'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.52%  '
$ws.Range("B49").Value = 'EOS'
$ws.Range("C49").Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range("D49").Value = '''1.225'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.75%  '
$ws.Range("B50").Value = 'NEARProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D50").Value = '''1.991'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.00%  '
$ws.Range("D51").Value = '''120.77'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.51%  '